$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet gained two new rows (26 and 27), each a duplicate of row 25
# ("Ridvan"), pushing the used range from A1:G25 to A1:G27.
# Using Rows.Copy()/Insert() (instead of just writing .Value) so the new
# rows inherit row 25's cell style (s="2"), matching the target diff
# instead of landing with no style applied.
$ws.Rows.Item(25).Copy() | Out-Null
$ws.Rows.Item(26).Insert()
$ws.Rows.Item(25).Copy() | Out-Null
$ws.Rows.Item(27).Insert()
$excel.CutCopyMode = $false

# Reflect the final selection recorded in the diff (cell F27 selected).
$ws.Range("F27").Select() | Out-Null
